$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New loading-percent values for columns B, C, E, F, G, H, J, O across rows 2-25
# (Case with 380 kV)
$startRow = 2

$colB = @(15.18689417704505, 14.34360736476312, 13.79850449933102, 13.56964594329898, 13.53124287543591, 13.79544504566274, 14.90187737512548, 16.84984806042555, 18.14082974772012, 18.69699775065682, 18.90309079419274, 18.85890637102642, 18.71404385649391, 18.62472222631012, 18.10385301961543, 17.77631558978197, 17.58499862921, 17.51972092211977, 17.81148569270758, 18.75671632944149, 19.34814572751013, 19.03490959644645, 17.79559466546606, 16.34720345870164)
for ($i = 0; $i -lt $colB.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 2).Value = $colB[$i]
}

$colC = @(11.25757293631382, 10.66888469173825, 10.2895047978297, 10.13052476259798, 10.10386593520015, 10.28737828433381, 11.05838213466143, 12.4236519468448, 13.3326482195531, 13.72506867989328, 13.87059425215952, 13.83939009151336, 13.73710300919195, 13.67404756869593, 13.30657437917377, 13.07570674465065, 12.94093493315004, 12.8949644374341, 13.10048856083172, 13.76723097230872, 14.18504535318349, 13.96370342982238, 13.08929106169666, 12.07055566409732)
for ($i = 0; $i -lt $colC.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 3).Value = $colC[$i]
}

$colE = @(26.22868936680642, 25.88623601065564, 25.67871153912892, 25.59492630203211, 25.58106378570472, 25.6775782914669, 26.11009814225543, 26.97591249488998, 27.61725133761335, 27.90898719165797, 28.01936326581394, 27.99559766607008, 27.91807040342106, 27.87056714491723, 27.59817718966807, 27.43100245401105, 27.33485242516778, 27.30230141089027, 27.4487987327267, 27.94084545306014, 28.26181259897508, 28.09059369361841, 27.44075314934927, 26.74037063787305)
for ($i = 0; $i -lt $colE.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 5).Value = $colE[$i]
}

$colF = @(37.41752663735916, 37.25262696296107, 37.16338758619354, 37.13006756944736, 37.1247194196983, 37.16292585830421, 37.35819359917984, 37.83503916218888, 38.24048868110859, 38.43635515409575, 38.51211533628821, 38.49572920283447, 38.44255644505355, 38.41019199576295, 38.22791379346613, 38.11898168602492, 38.05740566070807, 38.03674392913594, 38.13046641237759, 38.45813186750284, 38.68151967233556, 38.56146644468441, 38.12527089127835, 37.69618828858803)
for ($i = 0; $i -lt $colF.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 6).Value = $colF[$i]
}

$colG = @(17.07343389154736, 17.30280405665209, 17.45608994382393, 17.52164384547619, 17.53271435032903, 17.45696157608901, 17.14991052354051, 16.64853469311452, 16.34434477612756, 16.22050853747436, 16.17575814343245, 16.18529975176587, 16.21678362781871, 16.2363491626495, 16.35273485647898, 16.42789527688719, 16.47249104395633, 16.48782333674291, 16.41975262079321, 16.2074774466506, 16.08126349947098, 16.14746304374765, 16.42342960521927, 16.77310466144584)
for ($i = 0; $i -lt $colG.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 7).Value = $colG[$i]
}

$colH = @(11.48091690534303, 11.56110944830529, 11.61313907350666, 11.6350437222983, 11.6387233910491, 11.61343164383371, 11.50798818809961, 11.32334561968047, 11.20116131349252, 11.14849892848347, 11.12897692526754, 11.13316264960053, 11.14688442058167, 11.15534411868611, 11.20466169712954, 11.23566442127434, 11.25377118226864, 11.25994901428419, 11.23233568470298, 11.14284260031383, 11.08680250946396, 11.11648800148096, 11.23383972532758, 11.37092921831264)
for ($i = 0; $i -lt $colH.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 8).Value = $colH[$i]
}

$colJ = @(7.1912830727934, 7.220656048891803, 7.239665583287085, 7.247657706113192, 7.248999641189167, 7.239772372750135, 7.201208974205171, 7.133292071125146, 7.088055076728181, 7.068480201354355, 7.061211475862541, 7.062770535353713, 7.067879318527193, 7.071027316096582, 7.089354496688853, 7.100854352527118, 7.107563262237865, 7.109851028182594, 7.099620396606672, 7.06637484459033, 7.045485239569565, 7.056557864806395, 7.100177964362302, 7.150844197364274)
for ($i = 0; $i -lt $colJ.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 10).Value = $colJ[$i]
}

$colO = @(15.89290304640388, 16.04406747924275, 16.14271073150358, 16.18437009530287, 16.19137575602757, 16.14326665208432, 15.94381255953992, 15.59910914520982, 15.37442414741553, 15.27847558138573, 15.24304868760294, 15.2506380711799, 15.27554277962749, 15.29091591078906, 15.38082114648022, 15.43758327022838, 15.47082033147828, 15.48217480195636, 15.43147984241297, 15.26820300226234, 15.16678052748925, 15.22042555079974, 15.43423732186263, 15.68736396300553)
for ($i = 0; $i -lt $colO.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 15).Value = $colO[$i]
}
